# verification.xlsx update
#  - fill in previously-blank "experimental result" B cells on interval_estimation
#    with simulation values / formulas
#  - fix "Global Premium response time" experimental value (C61) and drop the
#    leading space on the two "Global ..." analytical strings
#  - fix a typo ("premium center" -> "premium centerù")
#  - add new sheet "valori_attesi_maggiori" summarising the two metrics whose
#    experimental value fell outside the analytical confidence interval

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interval_estimation")

# --- Response time ------------------------------------------------------
$ws.Range("B7").Value = 232.945742
$ws.Range("B7").Font.Bold = $true

$ws.Range("B8").Value = 144.375297
$ws.Range("E8").Formula = "=143.158816+1.413922"

# --- Waiting time --------------------------------------------------------
$ws.Range("B11").Value = 85.693088
$ws.Range("B11").Font.Bold = $true

$ws.Range("B12").Value = 25.183851

# --- Global response time / Global Premium response time -----------------
$ws.Range("C29").Value = "165.793473 +/- 1.765061"
$ws.Range("C30").Value = "141.137177 +/- 1.397107"

$ws.Range("B29").Formula = "=TRUNC(B6 + B7*B39/B2 + B8*B40/B2 + B9*B41/B2,6)"
$ws.Range("B30").Formula = "=TRUNC(B6 + B8*B40/B2 + 148.418018*B37, 6)"

# --- Rho in the ... centers ------------------------------------------------
$ws.Range("B32").Formula = "=TRUNC(B14/B18,6)"
$ws.Range("B33").Formula = "=TRUNC(B15/(B19*50),6)"

$ws.Range("B34").Formula = "=TRUNC(B16/(B20*95),6)"
$ws.Range("B34").Font.Bold = $true

$ws.Range("B35").Formula = "=TRUNC(B17/(B21),6)"
$ws.Range("B35").Font.Bold = $true

# --- typo fix --------------------------------------------------------------
$ws.Range("A24").Value = "Average number of jobs in premium centerù"

# --- cosmetic: selection + column width ------------------------------------
$ws.Range("B31").Select()
$ws.Columns("B").ColumnWidth = 14.88

# --- new summary sheet -------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "valori_attesi_maggiori"

$newSheet.Range("A1").Value = "Response time normal center"
$newSheet.Range("B1").Value = 232.945742
$newSheet.Range("B1").Font.Bold = $true
$newSheet.Range("C1").Value = " 224.062908 +/- 4.812282"

$newSheet.Range("A2").Value = "Waiting time normal center"
$newSheet.Range("B2").Value = 85.693088
$newSheet.Range("B2").Font.Bold = $true
$newSheet.Range("C2").Value = " 76.828723 +/- 4.738207"

$newSheet.Range("B3").Font.Bold = $true

$newSheet.Columns("A").ColumnWidth = 32.09
$newSheet.Columns("C").ColumnWidth = 32.51
